# ============================================================================
# [ADDITIONAL SCRAPING] add a "Player Info" sheet and an "ODI Batting Extra"
# sheet, and replace the MATCH_CARD_LINK url columns with plain MATCH_CODE
# numbers on the existing "ODI Batting" / "ODI Bowling" sheets.
#
# NOTE: worksheet object handles in this host are positional, not stable
# identities -- once the tab order/count changes (Add/Move/rename), any
# previously-grabbed $sheet variable can silently start pointing at whatever
# now sits at that same index. So every step below re-fetches the sheet it
# needs by name right before touching it, instead of caching references
# across structural operations.
# ============================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($cell, $text) {
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin)
}

# Digit-ish values (match codes, counts, percentages, ...) need the cell
# pre-formatted as Text, otherwise Excel "helpfully" reinterprets "4253" as
# a number or "4.66%" as a percentage instead of keeping the literal text.
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ----------------------------------------------------------------------
# 1. Create the sheet skeleton / final tab order first:
#    Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ----------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("ODI Batting")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "Player Info"

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "ODI Batting Extra"

# ----------------------------------------------------------------------
# 2. Populate "Player Info".
# ----------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")

Set-HeaderCell $playerInfo.Cells.Item(1,1) "ID"
Set-HeaderCell $playerInfo.Cells.Item(1,2) "NAME"
Set-HeaderCell $playerInfo.Cells.Item(1,3) "BATTING_HAND"
Set-HeaderCell $playerInfo.Cells.Item(1,4) "BOWL_STYLE"

Set-TextCell $playerInfo.Cells.Item(2,1) "4587"
$playerInfo.Cells.Item(2,2).Value = "Nicholas Pooran"
$playerInfo.Cells.Item(2,3).Value = "Left Handed"
$playerInfo.Cells.Item(2,4).Value = "Does Not Bowl | Unknown"

# ----------------------------------------------------------------------
# 3. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#    scorecard urls with the bare match code extracted from them. Also
#    drop the stray empty INNING_NUMBER cells on the "did not bat" rows.
# ----------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

Set-HeaderCell $battingSheet.Cells.Item(1,4) "MATCH_CODE"

$battingRows = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingRows; $r++) {
    $link = $battingSheet.Cells.Item($r, 4).Value2
    if ($link) {
        $code = $link -replace '^.*MatchCode=', ''
        Set-TextCell $battingSheet.Cells.Item($r, 4) $code
    }

    $inning = $battingSheet.Cells.Item($r, 2).Value2
    if (-not $inning) {
        $battingSheet.Cells.Item($r, 2).ClearContents()
    }
}

# ----------------------------------------------------------------------
# 4. "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE treatment.
# ----------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

Set-HeaderCell $bowlingSheet.Cells.Item(1,2) "MATCH_CODE"

$bowlingRows = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingRows; $r++) {
    $link = $bowlingSheet.Cells.Item($r, 2).Value2
    if ($link) {
        $code = $link -replace '^.*MatchCode=', ''
        Set-TextCell $bowlingSheet.Cells.Item($r, 2) $code
    }
}

# ----------------------------------------------------------------------
# 5. Populate "ODI Batting Extra" with per-innings batting-position /
#    boundary-count / match-award detail.
# ----------------------------------------------------------------------
$extraSheet = $wb.Worksheets.Item("ODI Batting Extra")

Set-HeaderCell $extraSheet.Cells.Item(1,1) "MATCH_CODE"
Set-HeaderCell $extraSheet.Cells.Item(1,2) "BATTING_POSITION"
Set-HeaderCell $extraSheet.Cells.Item(1,3) "NUM_4"
Set-HeaderCell $extraSheet.Cells.Item(1,4) "NUM_6"
Set-HeaderCell $extraSheet.Cells.Item(1,5) "PERCENT_RUNS_OF_TOTAL"
Set-HeaderCell $extraSheet.Cells.Item(1,6) "MAN_OF_MATCH"

# MATCH_CODE, BATTING_POSITION (numeric or blank), NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4533", $null, $null, $null, $null, "NO"),
    @("4535", 5, "0", "1", "4.66%", "NO"),
    @("4536", 5, "2", "1", "20.12%", "NO"),
    @("4577", 4, "0", "0", "2.81%", "NO"),
    @("4580", 4, "0", "1", "4.61%", "NO"),
    @("4583", 4, "1", "0", "2.27%", "NO"),
    @("4586", $null, $null, $null, $null, "NO"),
    @("4590", $null, $null, $null, $null, "NO"),
    @("4592", 5, "1", "0", "5.09%", "NO"),
    @("4606", $null, $null, $null, $null, "NO"),
    @("4611", 5, "0", "0", $null, "NO"),
    @("4616", 5, "4", "2", "41.01%", "NO"),
    @("4621", 5, "0", "2", "8.20%", "NO"),
    @("4623", $null, $null, $null, $null, "NO"),
    @("4624", 5, "5", "1", "30.66%", "NO"),
    @("4636", 5, "2", "0", "14.51%", "NO"),
    @("4639", 6, "0", "0", "1.24%", "NO"),
    @("4642", $null, $null, $null, $null, "NO"),
    @("4727", $null, $null, $null, $null, "NO"),
    @("4731", 5, "1", "3", "15.00%", "NO")
)

$r = 2
foreach ($row in $extraRows) {
    Set-TextCell $extraSheet.Cells.Item($r, 1) $row[0]

    if ($null -ne $row[1]) {
        $extraSheet.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        Set-TextCell $extraSheet.Cells.Item($r, 3) $row[2]
    }
    if ($null -ne $row[3]) {
        Set-TextCell $extraSheet.Cells.Item($r, 4) $row[3]
    }
    if ($null -ne $row[4]) {
        Set-TextCell $extraSheet.Cells.Item($r, 5) $row[4]
    }
    $extraSheet.Cells.Item($r, 6).Value = $row[5]

    $r++
}
